$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.604199999999991
$ws.Range("C6").Value = -11.7169
$ws.Range("C7").Value = -11.83779999999999
$ws.Range("B8").Value = 4.937699999999997
$ws.Range("C8").Value = -11.26179999999999
$ws.Range("A12").Value = -22.65530000000001
$ws.Range("B12").Value = 5.632199999999998
$ws.Range("B14").Value = 8.594100000000005
$ws.Range("C19").Value = -13.22559999999999
$ws.Range("C21").Value = -13.0606
$ws.Range("B22").Value = 4.814300000000006
$ws.Range("C24").Value = -11.55119999999999
